# Regenerate LinkML artifacts (#41)
#
# Adds four new worksheets that mirror existing "Acquisition*" sheets but
# with an extra "technique" column, matching the regenerated LinkML schema
# export:
#   - AcquisitionSpa      : copy of "Acquisition" + "technique" column A
#                            (inserted right after "Acquisition")
#   - AcquisitionTomo     : existing sheet gets a new "technique" column
#                            inserted right after "tilt_angle" (column C)
#   - AcquisitionSubTomo  : copy of the (now updated) "AcquisitionTomo" row
#   - AcquisitionEnvTomo  : copy of the (now updated) "AcquisitionTomo" row
#   - AcquisitionCelTomo  : copy of the (now updated) "AcquisitionTomo" row
#                            (each inserted, in order, right after
#                            "AcquisitionTomo")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) AcquisitionSpa = Acquisition's header row + a new "technique" column
#    at the front, placed right after the "Acquisition" tab.
# ---------------------------------------------------------------------
$acquisition = $wb.Worksheets.Item("Acquisition")
$acquisition.Range("A1:AA1").Copy()

$acquisitionSpa = $wb.Worksheets.Add($null, $acquisition)
$acquisitionSpa.Name = "AcquisitionSpa"
$acquisitionSpa.Range("B1").PasteSpecial()
$acquisitionSpa.Range("A1").Value = "technique"

# ---------------------------------------------------------------------
# 2) AcquisitionTomo: insert a "technique" column right after tilt_angle
#    (i.e. before the old column C / screen_current), shifting everything
#    from screen_current onward one column to the right.
# ---------------------------------------------------------------------
$acquisitionTomo = $wb.Worksheets.Item("AcquisitionTomo")
$acquisitionTomo.Range("C1").EntireColumn.Insert()
$acquisitionTomo.Range("C1").Value = "technique"

# ---------------------------------------------------------------------
# 3) AcquisitionSubTomo, AcquisitionEnvTomo, AcquisitionCelTomo: each is a
#    copy of the now-updated AcquisitionTomo header row, inserted in turn
#    right after AcquisitionTomo (so final order is AcquisitionTomo ->
#    AcquisitionSubTomo -> AcquisitionEnvTomo -> AcquisitionCelTomo).
# ---------------------------------------------------------------------
$acquisitionTomo.Range("A1:AD1").Copy()
$acquisitionSubTomo = $wb.Worksheets.Add($null, $acquisitionTomo)
$acquisitionSubTomo.Name = "AcquisitionSubTomo"
$acquisitionSubTomo.Range("A1").PasteSpecial()

$acquisitionSubTomo.Range("A1:AD1").Copy()
$acquisitionEnvTomo = $wb.Worksheets.Add($null, $acquisitionSubTomo)
$acquisitionEnvTomo.Name = "AcquisitionEnvTomo"
$acquisitionEnvTomo.Range("A1").PasteSpecial()

$acquisitionEnvTomo.Range("A1:AD1").Copy()
$acquisitionCelTomo = $wb.Worksheets.Add($null, $acquisitionEnvTomo)
$acquisitionCelTomo.Name = "AcquisitionCelTomo"
$acquisitionCelTomo.Range("A1").PasteSpecial()

# Restore the originally-active sheet/selection (sheet-adding operations
# above shift the active tab as a side effect).
$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Activate()
$null = $firstSheet.Range("A1").Select()
